$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.516.32'
$ws.Range('E2').Value = '  -1.04%  '

$ws.Range('D3').Value = '2.928.22'
$ws.Range('E3').Value = '  -2.69%  '

$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '374.67'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.63%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.89%  '

$ws.Range('E7').Value = '  -2.81%  '

$ws.Range('E8').Value = '  -0.23%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.585'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.49%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.00'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.05%  '

$ws.Range('E11').Value = '  -0.59%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0837'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.50%  '

$ws.Range('E13').Value = '  -3.71%  '

$ws.Range('D14').Value = '3.389.95'
$ws.Range('E14').Value = '  -2.54%  '

$ws.Range('E15').Value = '  -3.75%  '

$ws.Range('D16').Value = '2.923.78'
$ws.Range('E16').Value = '  -3.11%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.929'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -8.63%  '

$ws.Range('D18').Value = '51.451.45'
$ws.Range('E18').Value = '  -1.29%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.43'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.92%  '

$ws.Range('E20').Value = '  -1.90%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.96'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.66%  '

$ws.Range('E22').Value = '  -2.86%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.28%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '262.19'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.69%  '

$ws.Range('E25').Value = '  +0.82%  '

$ws.Range('E26').Value = '  -5.93%  '

$ws.Range('E27').Value = '  -4.85%  '

$ws.Range('E28').Value = '  -0.03%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '25.74'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.74%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.30'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.07%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.95'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.92%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.101'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.12%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.82'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.66%  '

$ws.Range('E34').Value = '  -3.42%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '51.03'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.35%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '34.01'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.97%  '

$ws.Range('E37').Value = '  +0.41%  '

$ws.Range('E38').Value = '  -3.58%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.00'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -10.15%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.98'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.77%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.57'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -9.52%  '

$ws.Range('E42').Value = '  -7.68%  '

$ws.Range('E43').Value = '  -2.34%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '121.81'
$ws.Range('D44').Style = 'Normal'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.75'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.77%  '

$ws.Range('E46').Value = '  -3.91%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.273'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +12.01%  '

$ws.Range('D48').Value = '2.023.02'
$ws.Range('E48').Value = '  -4.78%  '

$ws.Range('E49').Value = '  -1.85%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.16'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.21%  '

$ws.Range('D51').Value = '3.207.05'
